# Update the numeric data table (columns A:E, rows 2-51) on Sheet1 to the
# new values produced by the author's re-run of the analysis. Row 47 keeps
# its original E47 value (unchanged in the source diff); every other A:E
# cell in rows 2-51 receives a new value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.033488205426709321
$ws.Range("B2").Value = 1.6573876655055337
$ws.Range("C2").Value = 1.6200488214049564
$ws.Range("D2").Value = 1.514164889895498
$ws.Range("E2").Value = -0.0000099999999999964517
$ws.Range("A3").Value = 0.090769099612936457
$ws.Range("B3").Value = 1.7882324643576384
$ws.Range("C3").Value = 1.700919232560256
$ws.Range("D3").Value = 1.5388309715970907
$ws.Range("E3").Value = -0.0000099999999999969363
$ws.Range("A4").Value = 0.16516032460820018
$ws.Range("B4").Value = 1.910973836650006
$ws.Range("C4").Value = 1.7733501337212731
$ws.Range("D4").Value = 1.5723238294518531
$ws.Range("E4").Value = -0.0000099999999999961231
$ws.Range("A5").Value = 0.22635885034944389
$ws.Range("B5").Value = 2.0120817236021336
$ws.Range("C5").Value = 1.832075349157313
$ws.Range("D5").Value = 1.6119020399253783
$ws.Range("E5").Value = 0.14648745553971695
$ws.Range("A6").Value = 0.24992084702376424
$ws.Range("B6").Value = 2.074014150281954
$ws.Range("C6").Value = 1.8737954754274537
$ws.Range("D6").Value = 1.6534921298574701
$ws.Range("E6").Value = 0.35642061867104624
$ws.Range("A7").Value = 0.23978613835948714
$ws.Range("B7").Value = 2.0965754683781084
$ws.Range("C7").Value = 1.8990045530898327
$ws.Range("D7").Value = 1.6933696303243275
$ws.Range("E7").Value = 0.48164652485099568
$ws.Range("A8").Value = 0.21030707432105028
$ws.Range("B8").Value = 2.0915092251004235
$ws.Range("C8").Value = 1.9121026729121597
$ws.Range("D8").Value = 1.7292021938541762
$ws.Range("E8").Value = 0.52915482804008873
$ws.Range("A9").Value = 0.17550581645697874
$ws.Range("B9").Value = 2.0733268137091172
$ws.Range("C9").Value = 1.918549084468727
$ws.Range("D9").Value = 1.760173409641121
$ws.Range("E9").Value = 0.52645840883100115
$ws.Range("A10").Value = 0.14445152085888557
$ws.Range("B10").Value = 2.0533338527440046
$ws.Range("C10").Value = 1.9227184010319971
$ws.Range("D10").Value = 1.7865578495203804
$ws.Range("E10").Value = 0.50250334972984656
$ws.Range("A11").Value = 0.12079314247462751
$ws.Range("B11").Value = 2.0376571570136939
$ws.Range("C11").Value = 1.9271154910336268
$ws.Range("D11").Value = 1.8091567871947787
$ws.Range("E11").Value = 0.47772419226448648
$ws.Range("A12").Value = 0.10434109574713313
$ws.Range("B12").Value = 2.0279533862063412
$ws.Range("C12").Value = 1.9325789896771899
$ws.Range("D12").Value = 1.8288484811058192
$ws.Range("E12").Value = 0.46193950463197841
$ws.Range("A13").Value = 0.093085048777225204
$ws.Range("B13").Value = 2.023195555525088
$ws.Range("C13").Value = 1.9389110926285533
$ws.Range("D13").Value = 1.8463397178035534
$ws.Range("E13").Value = 0.45675491150406788
$ws.Range("A14").Value = 0.084755400293326275
$ws.Range("B14").Value = 2.0213980027328438
$ws.Range("C14").Value = 1.9455106934577913
$ws.Range("D14").Value = 1.8620949634471895
$ws.Range("E14").Value = 0.45936221004051464
$ws.Range("A15").Value = 0.07765697948701196
$ws.Range("B15").Value = 2.0207517225958886
$ws.Range("C15").Value = 1.9518034150084116
$ws.Range("D15").Value = 1.8763740717705726
$ws.Range("E15").Value = 0.46577836472872614
$ws.Range("A16").Value = 0.070883866381532878
$ws.Range("B16").Value = 2.0201093588094317
$ws.Range("C16").Value = 1.9574354104869285
$ws.Range("D16").Value = 1.88931024760407
$ws.Range("E16").Value = 0.47276784843905889
$ws.Range("A17").Value = 0.064168919136775512
$ws.Range("B17").Value = 2.0190023817166916
$ws.Range("C17").Value = 1.9622918255522432
$ws.Range("D17").Value = 1.9009825396742059
$ws.Range("E17").Value = 0.47851381104397489
$ws.Range("A18").Value = 0.057607755483568272
$ws.Range("B18").Value = 2.017429602387669
$ws.Range("C18").Value = 1.9664259214045989
$ws.Range("D18").Value = 1.9114627753184177
$ws.Range("E18").Value = 0.48247536774889072
$ws.Range("A19").Value = 0.051412300290928678
$ws.Range("B19").Value = 2.0156053618207128
$ws.Range("C19").Value = 1.9699688709923135
$ws.Range("D19").Value = 1.9208356081036242
$ws.Range("E19").Value = 0.48490373264472236
$ws.Range("A20").Value = 0.045758230688951816
$ws.Range("B20").Value = 2.0137702701887856
$ws.Range("C20").Value = 1.9730598506052397
$ws.Range("D20").Value = 1.9291997276912887
$ws.Range("E20").Value = 0.48635325968919124
$ws.Range("A21").Value = 0.040726473588489118
$ws.Range("B21").Value = 2.0120931817988272
$ws.Range("C21").Value = 1.9758086939689155
$ws.Range("D21").Value = 1.9366601385609672
$ws.Range("E21").Value = 0.48734911982212709
$ws.Range("A22").Value = 0.036308281639759334
$ws.Range("B22").Value = 2.0106484309397934
$ws.Range("C22").Value = 1.9782863374358695
$ws.Range("D22").Value = 1.9433190848750617
$ws.Range("E22").Value = 0.48823772933549975
$ws.Range("A23").Value = 0.03243872265461261
$ws.Range("B23").Value = 2.0094362307319713
$ws.Range("C23").Value = 1.980531749078104
$ws.Range("D23").Value = 1.9492696280021835
$ws.Range("E23").Value = 0.48917292536110113
$ws.Range("A24").Value = 0.029032618640441179
$ws.Range("B24").Value = 2.0084169535265208
$ws.Range("C24").Value = 1.9825645397884815
$ws.Range("D24").Value = 1.9545928872993739
$ws.Range("E24").Value = 0.49017137192825599
$ws.Range("A25").Value = 0.026009782320742074
$ws.Range("B25").Value = 2.0075407477544758
$ws.Range("C25").Value = 1.9843963124071922
$ws.Range("D25").Value = 1.9593581947403331
$ws.Range("E25").Value = 0.49118293094921134
$ws.Range("A26").Value = 0.023306797004093393
$ws.Range("B26").Value = 2.0067653126865199
$ws.Range("C26").Value = 1.986037701538607
$ws.Range("D26").Value = 1.9636248353554899
$ws.Range("E26").Value = 0.49214535175199831
$ws.Range("A27").Value = 0.020878395202665316
$ws.Range("B27").Value = 2.0060619734825962
$ws.Range("C27").Value = 1.9875012674787538
$ws.Range("D27").Value = 1.9674441777869294
$ws.Range("E27").Value = 0.49301357682048474
$ws.Range("A28").Value = 0.018693334352696257
$ws.Range("B28").Value = 2.0054144164445882
$ws.Range("C28").Value = 1.9888011959946121
$ws.Range("D28").Value = 1.9708614992661189
$ws.Range("E28").Value = 0.49376758108430108
$ws.Range("A29").Value = 0.016728821100437885
$ws.Range("B29").Value = 2.0048135627889287
$ws.Range("C29").Value = 1.9899520832282225
$ws.Range("D29").Value = 1.9739171849831716
$ws.Range("E29").Value = 0.49440694041077454
$ws.Range("A30").Value = 0.014966056864819664
$ws.Range("B30").Value = 2.0042532072741222
$ws.Range("C30").Value = 1.990966599818639
$ws.Range("D30").Value = 1.9766474269893572
$ws.Range("E30").Value = 0.49494139575592377
$ws.Range("A31").Value = 0.013387482472846978
$ws.Range("B31").Value = 2.0037250077073105
$ws.Range("C31").Value = 1.9918550883554424
$ws.Range("D31").Value = 1.979084409253973
$ws.Range("E31").Value = 0.49538095114814523
$ws.Range("A32").Value = 0.011976002735080865
$ws.Range("B32").Value = 2.0032183617393464
$ws.Range("C32").Value = 1.9926224974891862
$ws.Range("D32").Value = 1.9812564649776565
$ws.Range("E32").Value = 0.4957310657626196
$ws.Range("A33").Value = 0.010714670833336293
$ws.Range("B33").Value = 2.0027146576746513
$ws.Range("C33").Value = 1.9932710837948331
$ws.Range("D33").Value = 1.9831877023203861
$ws.Range("E33").Value = 0.49598652823022837
$ws.Range("A34").Value = 0.009587864139436162
$ws.Range("B34").Value = 2.0021936421755511
$ws.Range("C34").Value = 1.9937923911141822
$ws.Range("D34").Value = 1.9848982369073511
$ws.Range("E34").Value = 0.49613447047209647
$ws.Range("A35").Value = 0.0085804369511346366
$ws.Range("B35").Value = 2.0016166230662598
$ws.Range("C35").Value = 1.9941781077946739
$ws.Range("D35").Value = 1.986402891661653
$ws.Range("E35").Value = 0.49614317803160873
$ws.Range("A36").Value = 0.0076804077896918625
$ws.Range("B36").Value = 2.0009500018463493
$ws.Range("C36").Value = 1.9943935497089904
$ws.Range("D36").Value = 1.9877121315782755
$ws.Range("E36").Value = 0.49597589324635877
$ws.Range("A37").Value = 0.006874470260309273
$ws.Range("B37").Value = 2.0001090101169865
$ws.Range("C37").Value = 1.9944171309182339
$ws.Range("D37").Value = 1.9888278506467596
$ws.Range("E37").Value = 0.4955529783278661
$ws.Range("A38").Value = 0.006156298317122339
$ws.Range("B38").Value = 1.9990430290566601
$ws.Range("C38").Value = 1.9941518930418045
$ws.Range("D38").Value = 1.9897472167036507
$ws.Range("E38").Value = 0.49480178344934
$ws.Range("A39").Value = 0.0055107146399101031
$ws.Range("B39").Value = 1.9975478229745172
$ws.Range("C39").Value = 1.9935676940040663
$ws.Range("D39").Value = 1.9904492712680286
$ws.Range("E39").Value = 0.49353153600036698
$ws.Range("A40").Value = 0.0049426817262087216
$ws.Range("B40").Value = 1.9955721997707918
$ws.Range("C40").Value = 1.9924035996646419
$ws.Range("D40").Value = 1.9909103348332773
$ws.Range("E40").Value = 0.49161748716734954
$ws.Range("A41").Value = 0.0044229891754760204
$ws.Range("B41").Value = 1.9925850351111016
$ws.Range("C41").Value = 1.990673149104369
$ws.Range("D41").Value = 1.9910610578582815
$ws.Range("E41").Value = 0.4885866608024535
$ws.Range("A42").Value = 0.0039891318981817546
$ws.Range("B42").Value = 1.9886619183397949
$ws.Range("C42").Value = 1.9876533032520141
$ws.Range("D42").Value = 1.9908451353016177
$ws.Range("E42").Value = 0.48428390890429407
$ws.Range("A43").Value = 0.0035576797523223132
$ws.Range("B43").Value = 1.9823303226981601
$ws.Range("C43").Value = 1.9836618023077848
$ws.Range("D43").Value = 1.9900788021673066
$ws.Range("E43").Value = 0.47747977505070704
$ws.Range("A44").Value = 0.0032737186644082117
$ws.Range("B44").Value = 1.974384509972543
$ws.Range("C44").Value = 1.9765983005094541
$ws.Range("D44").Value = 1.9886663158697779
$ws.Range("E44").Value = 0.46823739427087629
$ws.Range("A45").Value = 0.0028613260560769026
$ws.Range("B45").Value = 1.9605032205647124
$ws.Range("C45").Value = 1.9681577209913692
$ws.Range("D45").Value = 1.9861316372923219
$ws.Range("E45").Value = 0.45318211860679208
$ws.Range("A46").Value = 0.0028331367179195003
$ws.Range("B46").Value = 1.9445694995580771
$ws.Range("C46").Value = 1.9519444635968419
$ws.Range("D46").Value = 1.9823910448962399
$ws.Range("E46").Value = 0.43382535613530637
$ws.Range("A47").Value = 0.002241428689335496
$ws.Range("B47").Value = 1.9132795095940298
$ws.Range("C47").Value = 1.9351115139434512
$ws.Range("D47").Value = 1.976171006719041
$ws.Range("A48").Value = 0.0028649030629038594
$ws.Range("B48").Value = 1.8824959689008873
$ws.Range("C48").Value = 1.8974025574855045
$ws.Range("D48").Value = 1.9677402533154071
$ws.Range("E48").Value = 0.36084714607413787
$ws.Range("A49").Value = 0.0014223327513893924
$ws.Range("B49").Value = 1.8097010795526909
$ws.Range("C49").Value = 1.8663405437848
$ws.Range("D49").Value = 1.9535167276767629
$ws.Range("E49").Value = 0.28540303313969961
$ws.Range("A50").Value = 0.0040652740308750419
$ws.Range("B50").Value = 1.7549837902459493
$ws.Range("C50").Value = 1.7758417722501996
$ws.Range("D50").Value = 1.9356487633079895
$ws.Range("E50").Value = 0.20769857954848889
$ws.Range("A51").Value = -0.00057933119613447649
$ws.Range("B51").Value = 1.5787205495064089
$ws.Range("C51").Value = 1.7271013074309058
$ws.Range("D51").Value = 1.9035252240658473
$ws.Range("E51").Value = 0.033161788153726883
